# Applies the cryptos.xlsx price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.751.36"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "1.564.72"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "206.22"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.489"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -2.06%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.788.00"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.568.62"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "61.53"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "26.787.25"
$ws.Range("E17").Value = "  -2.19%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "214.00"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.83%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.34"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("E24").Value = "  -0.89%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "152.36"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").Value = "1.385.26"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("E35").Value = "  +1.15%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.924"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("E38").Value = "  -2.29%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.523"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.77%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.816"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  +0.01%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.990"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.33"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("E45").Value = "  +1.23%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "63.20"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").Value = "1.699.99"
$ws.Range("E47").Value = "  +0.28%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "85.36"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "0.0₇0984"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  -0.71%  "
